{"js": "// Office.js (Word JavaScript API) \u2014 body of `async (context) => { ... }`.\n//\n// Implements the commit's changes:\n//  1. Paragraph \"git add .\": split its single run into two runs\n//     (\"git \" and \"add .\"), wrap the second run with\n//     <w:proofErr w:type=\"gramStart\"/> ... <w:proofErr w:type=\"gramEnd\"/>,\n//     and wrap the \"git add .\" .. \"git push origin main\" paragraphs in an\n//     \"OLE_LINK1\" bookmark (id 0).\n//  2. Table cell \"window.location.href='login.html'\": split the single\n//     run into three runs (\"window.location\", \".href\", \"='login.html'\")\n//     with proofErr spellStart/gramStart before the first run, gramEnd\n//     after the first run, and spellEnd after the second run.\n\nfunction wrapOoxml(bodyXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    bodyXml +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nasync function findParagraphByText(body, text) {\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n  paragraphs.items.forEach((p) => p.load(\"text\"));\n  await context.sync();\n  for (const p of paragraphs.items) {\n    if (p.text === text) {\n      return p;\n    }\n  }\n  return null;\n}\n\n// --- Step 1: split the \"git add .\" run and tag it with proofErr ------\nconst gitAddPara = await findParagraphByText(context.document.body, \"git add .\");\nif (!gitAddPara) {\n  throw new Error('Could not find the \"git add .\" paragraph.');\n}\n\nconst gitAddXml =\n  \"<w:p>\" +\n  '<w:pPr><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\">git </w:t></w:r>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>add .</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  \"</w:p>\";\n\ngitAddPara.getRange(\"Whole\").insertOoxml(wrapOoxml(gitAddXml), Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Step 2: re-resolve the paragraphs and wrap OLE_LINK1 bookmark ---\nconst gitAddPara2 = await findParagraphByText(context.document.body, \"git add .\");\nconst gitPushPara2 = await findParagraphByText(context.document.body, \"git push origin main\");\nif (!gitAddPara2 || !gitPushPara2) {\n  throw new Error(\"Could not re-locate the git paragraphs for bookmarking.\");\n}\n\nconst bookmarkRange = gitAddPara2.getRange(\"Start\").expandTo(gitPushPara2.getRange(\"End\"));\nbookmarkRange.insertBookmark(\"OLE_LINK1\");\nawait context.sync();\n\n// --- Step 3: split the \"window.location.href='login.html'\" run ------\nconst cellPara = await findParagraphByText(\n  context.document.body,\n  \"window.location.href='login.html'\"\n);\nif (!cellPara) {\n  throw new Error(\"Could not find the window.location.href paragraph.\");\n}\n\nconst cellXml =\n  \"<w:p>\" +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  \"<w:r><w:t>window.location</w:t></w:r>\" +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  \"<w:r><w:t>.href</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  \"<w:r><w:t>='login.html'</w:t></w:r>\" +\n  \"</w:p>\";\n\ncellPara.getRange(\"Whole\").insertOoxml(wrapOoxml(cellXml), Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word / $d / $app resolve against the open document.\n#\n# Implements the commit's changes:\n#  1. Paragraph \"git add .\": split its single run into two runs\n#     (\"git \" and \"add .\"), wrap the second run with\n#     <w:proofErr w:type=\"gramStart\"/> ... <w:proofErr w:type=\"gramEnd\"/>,\n#     and wrap the \"git add .\" .. \"git push origin main\" paragraphs in an\n#     \"OLE_LINK1\" bookmark (id 0).\n#  2. Table cell \"window.location.href='login.html'\": split the single\n#     run into three runs (\"window.location\", \".href\", \"='login.html'\")\n#     with proofErr spellStart/gramStart before the first run, gramEnd\n#     after the first run, and spellEnd after the second run.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphByText($doc, [string]$text) {\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        # Paragraph.Range.Text keeps its trailing paragraph mark (CR, chr 13)\n        # and, for the last paragraph of a table cell, a trailing end-of-cell\n        # mark (BEL, chr 7) as well \u2014 strip both before comparing.\n        $raw = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($raw -eq $text) {\n            return $p\n        }\n    }\n    return $null\n}\n\nfunction Wrap-Ooxml([string]$bodyXml) {\n    return '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $bodyXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\n# --- Step 1: split the \"git add .\" run and tag it with proofErr --------\n$gitAddPara = Find-ParagraphByText $d \"git add .\"\nif ($null -eq $gitAddPara) {\n    throw \"Could not find the 'git add .' paragraph.\"\n}\n\n$gitAddXml = '<w:p>' +\n    '<w:pPr><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\">git </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>add .</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '</w:p>'\n\n$gitAddPara.Range.InsertXML((Wrap-Ooxml $gitAddXml))\n\n# --- Step 2: re-resolve the paragraphs and wrap the OLE_LINK1 bookmark --\n$gitAddPara2 = Find-ParagraphByText $d \"git add .\"\n$gitPushPara2 = Find-ParagraphByText $d \"git push origin main\"\nif ($null -eq $gitAddPara2 -or $null -eq $gitPushPara2) {\n    throw \"Could not re-locate the git paragraphs for bookmarking.\"\n}\n\n$bmRange = $d.Range($gitAddPara2.Range.Start, $gitPushPara2.Range.End)\n$d.Bookmarks.Add(\"OLE_LINK1\", $bmRange)\n\n# --- Step 3: split the \"window.location.href='login.html'\" run --------\n$cellPara = Find-ParagraphByText $d \"window.location.href='login.html'\"\nif ($null -eq $cellPara) {\n    throw \"Could not find the window.location.href paragraph.\"\n}\n\n$quote = [char]39\n$cellXml = '<w:p>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>window.location</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t>.href</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>=' + $quote + \"login.html\" + $quote + '</w:t></w:r>' +\n    '</w:p>'\n\n$cellPara.Range.InsertXML((Wrap-Ooxml $cellXml))\n\n\"done\"\n"}
